$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "Circulation: Cardiovascular Interventions"
$ws.Range("G3").Value = "https://openalex.org/S177226689"
$ws.Range("H3").Value = "Lippincott Williams & Wilkins"
$ws.Range("I3").Value = "1941-7640"

# V3 holds the literal text "FALSE" (not a boolean). A bare Value = "FALSE"
# gets auto-typed to a boolean by the engine, so force text entry with a
# leading apostrophe, then clear the resulting "quote prefix" cell style so
# the cell keeps its original (unstyled) formatting.
$v3 = $ws.Range("V3")
$v3.Value = "'FALSE"
$v3.Style = "Normal"
